# Edit workbook per diff:
# - Delete the "theta_threshold_range" row (row 5), shifting
#   "pie_threshold_range" up to row 5.
# - Update the Min/Max numeric values for rows 2-5.
# - Move the active selection to C5 (previously D6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the theta_threshold_range row entirely; rows below shift up.
$ws.Rows.Item(5).Delete()

# Update numeric values to match the new data set.
$ws.Range("B2").Value = 5.2
$ws.Range("C2").Value = 9.5

$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 7.5

$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1.6

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 20

# Cells C4 and B5 (old B6) carried a one-off "Times New Roman" font style;
# normalize them back to the sheet's regular body style.
$ws.Range("C4").Font.Name = $ws.Range("A4").Font.Name
$ws.Range("C4").Font.Size = $ws.Range("A4").Font.Size
$ws.Range("B5").Font.Name = $ws.Range("A5").Font.Name
$ws.Range("B5").Font.Size = $ws.Range("A5").Font.Size

# Move selection to C5 to match the saved view state.
$ws.Range("C5").Select()
